# Generate Report for Handback
#
# The "Ready for handoff" status is now "Handed back: in sync with en-US"
# for both locales, and each locale sheet records the handback package
# (target file link, handback xliff file name, and handback datetime)
# for both rows (a.md, b.md).

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

$githubBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a4b7eb33382ab9ea8df9131743a266cd36cc722d/e2e/"

# ---- Overview sheet: refresh the per-locale status text ----------------
$overview = $wb.Worksheets.Item("Overview")
for ($r = 2; $r -le 3; $r++) {
    if ($overview.Cells.Item($r, 5).Value2 -eq $oldStatus) {
        $overview.Cells.Item($r, 5).Value = $newStatus
    }
    if ($overview.Cells.Item($r, 6).Value2 -eq $oldStatus) {
        $overview.Cells.Item($r, 6).Value = $newStatus
    }
}

# ---- Per-locale detail sheets -------------------------------------------
# columns: C=Status(3), I=Latest Target File(9), J=Latest Handback File(10),
#          K=Latest Handback DateTime(11)
$locales = @(
    @{ Sheet = "zh-cn"; XlfFile = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"; HandbackTime = "2016-08-22 20:36:39" },
    @{ Sheet = "de-de"; XlfFile = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"; HandbackTime = "2016-08-22 20:36:47" }
)

foreach ($locale in $locales) {
    $ws = $wb.Worksheets.Item($locale.Sheet)

    for ($r = 2; $r -le 3; $r++) {
        # Status text
        if ($ws.Cells.Item($r, 3).Value2 -eq $oldStatus) {
            $ws.Cells.Item($r, 3).Value = $newStatus
        }

        # Latest Target File -> hyperlink to a.md (display text always "a.md")
        $targetCell = $ws.Cells.Item($r, 9)
        $ws.Hyperlinks.Add($targetCell, ($githubBase + "a.md"), [Type]::Missing, [Type]::Missing, "a.md") | Out-Null

        # Latest Handback File
        $ws.Cells.Item($r, 10).Value = $locale.XlfFile

        # Latest Handback DateTime
        $ws.Cells.Item($r, 11).Value = $locale.HandbackTime
    }

    # Widen columns that now hold the longer status / file-name strings,
    # matching the new "Status" (C) and "Latest Handback File" (J) widths.
    $ws.Columns.Item(3).ColumnWidth = 29.9777047293527
    $ws.Columns.Item(10).ColumnWidth = 40
}

# Overview: widen the zh-cn / de-de status columns to fit the longer text.
$overview.Columns.Item(5).ColumnWidth = 29.9777047293527
$overview.Columns.Item(6).ColumnWidth = 29.9777047293527
